$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the commit diff (cryptos list refresh).
$ws.Range('D2').Value = '30.254.78'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '1.864.25'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '237.28'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = '0.4685'
$ws.Range('E7').Value = '  +0.84%  '
$ws.Range('D8').Value = '0.2867'
$ws.Range('E8').Value = '  +1.82%  '
$ws.Range('D9').Value = '0.06555'
$ws.Range('E9').Value = '  +0.31%  '
$ws.Range('D10').Value = '22.21'
$ws.Range('E10').Value = '  +11.83%  '
$ws.Range('D11').Value = '0.07905'
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').Value = '98.06'
$ws.Range('E12').Value = '  +2.08%  '
$ws.Range('D13').Value = '1.869.68'
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('E14').Value = '  +1.55%  '
$ws.Range('D15').Value = '0.6820'
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('D16').Value = '278.33'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('D17').Value = '30.254.07'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').Value = '13.64'
$ws.Range('E18').Value = '  +8.38%  '
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').Value = '0.000007344'
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.116.88'
$ws.Range('E21').Value = '  +1.05%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '5.370'
$ws.Range('E22').Value = '  -2.33%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = '6.203'
$ws.Range('E24').Value = '  +1.30%  '
$ws.Range('D25').Value = '168.15'
$ws.Range('E25').Value = '  +1.70%  '
$ws.Range('D26').Value = '9.246'
$ws.Range('E26').Value = '  -0.65%  '
$ws.Range('D27').Value = '19.10'
$ws.Range('E27').Value = '  +1.57%  '
$ws.Range('D28').Value = '1.954'
$ws.Range('E28').Value = '  +2.88%  '
$ws.Range('E29').Value = '  +2.97%  '
$ws.Range('D30').Value = '0.09862'
$ws.Range('E30').Value = '  +3.05%  '
$ws.Range('D31').Value = '4.393'
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').Value = '1.483'
$ws.Range('E32').Value = '  +1.12%  '
$ws.Range('D33').Value = '4.076'
$ws.Range('E33').Value = '  -0.46%  '
$ws.Range('E34').Value = '  +2.45%  '
$ws.Range('D35').Value = '1.139'
$ws.Range('E35').Value = '  +4.60%  '
$ws.Range('D36').Value = '0.7050'
$ws.Range('E36').Value = '  +1.11%  '
$ws.Range('D37').Value = '2.707'
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('D38').Value = '0.01883'
$ws.Range('E38').Value = '  +1.89%  '
$ws.Range('D39').Value = '2.631'
$ws.Range('E39').Value = '  +4.45%  '
$ws.Range('D40').Value = '76.07'
$ws.Range('E40').Value = '  +4.48%  '
$ws.Range('D41').Value = '6.293'
$ws.Range('E41').Value = '  +0.44%  '
$ws.Range('D42').Value = '1.961'
$ws.Range('E42').Value = '  +2.48%  '
$ws.Range('D43').Value = '0.8558'
$ws.Range('E43').Value = '  +0.46%  '
$ws.Range('D44').Value = '0.4187'
$ws.Range('E44').Value = '  +1.10%  '
$ws.Range('D45').Value = '0.9999'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').Value = '103.49'
$ws.Range('E46').Value = '  +0.23%  '
$ws.Range('D47').Value = '7.233'
$ws.Range('E47').Value = '  +1.15%  '
$ws.Range('D48').Value = '945.96'
$ws.Range('E48').Value = '  -4.41%  '
$ws.Range('D49').Value = '9.236'
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('D50').Value = '34.30'
$ws.Range('E50').Value = '  +0.77%  '
$ws.Range('E51').Value = '  +0.13%  '
